# RPA datasets push 2024-04-03
# Insert a new data row (row 2) for "아이엠비디엑스" (IMBDX) above the existing
# rows, shifting every other existing row down by one. All other rows keep
# their original values - only their row position changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 2 (pushes rows 2..16 down to 3..17)
$ws.Rows.Item(2).Insert()
# The inserted row inherits formatting from the row above (the bold/bordered
# header row) - strip that back to the plain default style used by the rest
# of the data rows.
$ws.Rows.Item(2).ClearFormats()

# Columns A-C hold plain text dates (e.g. "2024-03-14"). Force the cells to
# text format before assigning so Excel does not silently convert them into
# date serial numbers.
$ws.Range("A2:C2").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "2024-03-14"
$ws.Cells.Item(2, 2).Value = "2024-03-20"
$ws.Cells.Item(2, 3).Value = "2024-04-03"
$ws.Cells.Item(2, 4).Value = "미래"
$ws.Cells.Item(2, 5).Value = "아이엠비디엑스"
$ws.Cells.Item(2, 6).Value = 2500000
$ws.Cells.Item(2, 7).Value = 2500000
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 9).Value = 7700
$ws.Cells.Item(2, 10).Value = 9900
$ws.Cells.Item(2, 11).Value = 13992625
$ws.Cells.Item(2, 12).Value = 0
$ws.Cells.Item(2, 13).Value = 13000
$ws.Cells.Item(2, 14).Value = "865.73 :1"
$ws.Cells.Item(2, 15).Value = "-"
$ws.Cells.Item(2, 16).Value = 1230336508
$ws.Cells.Item(2, 17).Value = 2624739502
$ws.Cells.Item(2, 18).Value = 2926965114
$ws.Cells.Item(2, 19).Value = -5277789009
$ws.Cells.Item(2, 20).Value = -8667658271
$ws.Cells.Item(2, 21).Value = -4923399541
$ws.Cells.Item(2, 22).Value = -9788525741
$ws.Cells.Item(2, 23).Value = -10436419054
$ws.Cells.Item(2, 24).Value = -7563224846
$ws.Cells.Item(2, 25).Value = "알파리퀴드ⓡ 100, 알파리퀴드ⓡ HRR, 알파리퀴드ⓡ 디텍트, 알파리퀴드ⓡ 스크리닝"

# Restore the default (non-text) style on A2:C2 now that the text values are
# safely stored, so the cells match the unstyled look of the rest of the
# table instead of staying tagged with an explicit text-number-format style.
$ws.Range("A2:C2").Style = "Normal"
